$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: "каша гречневая рассыпчатая" -> "Каша рисовая молочная" ---
$ws.Range("C4").Value = "54-2о-2020"
$ws.Range("D4").Value = "Каша рисовая молочная"
$ws.Range("E4").Value = 200
$ws.Range("F4").Value = 16.42
$ws.Range("G4").Value = 215
$ws.Range("H4").Value = 13
$ws.Range("I4").Value = 15
$ws.Range("J4").Value = 7

# --- Row 5: "Кофейный напиток" -> "Кофейный напиток с молоком" ---
$ws.Range("C5").Value = "54-9ги2020"
$ws.Range("D5").Value = "Кофейный напиток с молоком"

# --- Row 6: "Хлеб в ассортименте" -> "хлеб в ассортименте" ---
$ws.Range("D6").Value = "хлеб в ассортименте"
$ws.Range("I6").Value = 0.3

# --- Row 7: "соус красный основной" -> "масло сливочное(порциями)" ---
$ws.Range("C7").Value = "53-19з2020"
$ws.Range("D7").Value = "масло сливочное(порциями)"
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 6
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = 0.1
$ws.Range("I7").Value = 8
$ws.Range("J7").Value = 0.1

# --- Row 8: remove the "54-2р-2020 / Котлеты мясная" dish entirely ---
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("G8").Value = ""
$ws.Range("H8").Value = ""
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = ""

# --- Row 9: fill in "Завтрак 2" / "фрукты" with "яблоко" nutrition data ---
$ws.Range("D9").Value = "яблоко"
$ws.Range("E9").Value = 150
$ws.Range("F9").Value = 22.5
$ws.Range("G9").Value = 61
$ws.Range("H9").Value = 1
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 14

# --- Row heights ---
$ws.Rows.Item(1).RowHeight = 13.8
$ws.Rows.Item(9).RowHeight = 13.8

# --- Selection ---
$ws.Range("J7").Select()

$wb.Save()
